$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) cells being updated to be stored as text so
# values such as '0.06780' or '1.0300' keep their exact digits instead of
# being reinterpreted/rounded as floating point numbers by Excel.
$priceCells = @('D2', 'D3', 'D4', 'D5', 'D6', 'D7', 'D8', 'D9', 'D10', 'D11', 'D12', 'D13', 'D14', 'D15', 'D16', 'D17', 'D18', 'D19', 'D20', 'D21', 'D22', 'D23', 'D24', 'D25', 'D26', 'D28', 'D29', 'D30', 'D31', 'D32', 'D33', 'D34', 'D35', 'D36', 'D37', 'D38', 'D39', 'D40', 'D41', 'D42', 'D43', 'D44', 'D45', 'D46', 'D47', 'D48', 'D49', 'D50', 'D51')
foreach ($pc in $priceCells) {
    $ws.Range($pc).NumberFormat = '@'
}

$ws.Range('D2').Value = '30.378.80'
$ws.Range('E2').Value = '  -2.76%  '
$ws.Range('D3').Value = '1.932.42'
$ws.Range('E3').Value = '  -2.64%  '
$ws.Range('D4').Value = '1.008'
$ws.Range('E4').Value = '  +1.01%  '
$ws.Range('D5').Value = '247.62'
$ws.Range('E5').Value = '  -2.31%  '
$ws.Range('D6').Value = '0.6912'
$ws.Range('E6').Value = '  -12.24%  '
$ws.Range('D7').Value = '1.007'
$ws.Range('E7').Value = '  +0.68%  '
$ws.Range('D8').Value = '0.3225'
$ws.Range('E8').Value = '  -5.05%  '
$ws.Range('D9').Value = '26.42'
$ws.Range('E9').Value = '  +3.79%  '
$ws.Range('D10').Value = '0.06780'
$ws.Range('E10').Value = '  -2.12%  '
$ws.Range('D11').Value = '0.7932'
$ws.Range('E11').Value = '  -6.85%  '
$ws.Range('D12').Value = '0.07995'
$ws.Range('E12').Value = '  -1.82%  '
$ws.Range('D13').Value = '1.945.91'
$ws.Range('E13').Value = '  -1.63%  '
$ws.Range('D14').Value = '5.379'
$ws.Range('E14').Value = '  -2.75%  '
$ws.Range('D15').Value = '94.07'
$ws.Range('E15').Value = '  -7.88%  '
$ws.Range('B16').Value = 'Avalanche'
$ws.Range('C16').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D16').Value = '14.44'
$ws.Range('E16').Value = '  +3.45%  '
$ws.Range('B17').Value = 'BitcoinCash'
$ws.Range('C17').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D17').Value = '260.48'
$ws.Range('E17').Value = '  -5.60%  '
$ws.Range('D18').Value = '30.396.76'
$ws.Range('E18').Value = '  -2.71%  '
$ws.Range('D19').Value = '5.867'
$ws.Range('E19').Value = '  +3.30%  '
$ws.Range('D20').Value = '0.000007794'
$ws.Range('E20').Value = '  -0.95%  '
$ws.Range('D21').Value = '2.205.21'
$ws.Range('E21').Value = '  -1.11%  '
$ws.Range('D22').Value = '1.007'
$ws.Range('E22').Value = '  +0.71%  '
$ws.Range('D23').Value = '1.008'
$ws.Range('E23').Value = '  +1.04%  '
$ws.Range('D24').Value = '6.815'
$ws.Range('E24').Value = '  +0.12%  '
$ws.Range('D25').Value = '9.607'
$ws.Range('E25').Value = '  -0.44%  '
$ws.Range('D26').Value = '158.85'
$ws.Range('E26').Value = '  -4.08%  '
$ws.Range('E27').Value = '  -4.30%  '
$ws.Range('B28').Value = 'Stellar'
$ws.Range('C28').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D28').Value = '0.1300'
$ws.Range('E28').Value = '  -15.84%  '
$ws.Range('B29').Value = 'LidoDAOToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D29').Value = '2.243'
$ws.Range('E29').Value = '  +1.15%  '
$ws.Range('D30').Value = '1.353'
$ws.Range('E30').Value = '  -0.14%  '
$ws.Range('D31').Value = '1.555'
$ws.Range('E31').Value = '  -0.74%  '
$ws.Range('D32').Value = '4.406'
$ws.Range('E32').Value = '  -3.21%  '
$ws.Range('D33').Value = '4.214'
$ws.Range('E33').Value = '  -2.85%  '
$ws.Range('D34').Value = '0.05077'
$ws.Range('E34').Value = '  -1.39%  '
$ws.Range('D35').Value = '1.198'
$ws.Range('E35').Value = '  -2.08%  '
$ws.Range('D36').Value = '0.7467'
$ws.Range('E36').Value = '  +0.72%  '
$ws.Range('D37').Value = '2.731'
$ws.Range('E37').Value = '  -2.55%  '
$ws.Range('D38').Value = '0.01926'
$ws.Range('E38').Value = '  -3.13%  '
$ws.Range('D39').Value = '2.785'
$ws.Range('E39').Value = '  -4.03%  '
$ws.Range('D40').Value = '79.91'
$ws.Range('E40').Value = '  +1.61%  '
$ws.Range('D41').Value = '6.522'
$ws.Range('E41').Value = '  -1.41%  '
$ws.Range('D42').Value = '2.041'
$ws.Range('E42').Value = '  -1.71%  '
$ws.Range('D43').Value = '0.4417'
$ws.Range('E43').Value = '  -5.18%  '
$ws.Range('D44').Value = '1.006'
$ws.Range('E44').Value = '  +0.50%  '
$ws.Range('D45').Value = '0.8396'
$ws.Range('E45').Value = '  -1.30%  '
$ws.Range('D46').Value = '101.77'
$ws.Range('E46').Value = '  -3.55%  '
$ws.Range('D47').Value = '9.710'
$ws.Range('E47').Value = '  -2.74%  '
$ws.Range('D48').Value = '7.280'
$ws.Range('E48').Value = '  -3.14%  '
$ws.Range('D49').Value = '35.92'
$ws.Range('E49').Value = '  -1.43%  '
$ws.Range('D50').Value = '1.488'
$ws.Range('E50').Value = '  +3.51%  '
$ws.Range('B51').Value = 'SynthetixNetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range('D51').Value = '2.805'
$ws.Range('E51').Value = '  +30.71%  '
